$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay as text (matches source
# workbook, which stores them as inline strings, e.g. "5.010", "1.005.03", etc.)
# Force the Text format on the whole updated range first so Excel does not
# silently coerce/round them into real numbers when we assign the new values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.616.78"
$ws.Range("E2").Value = "  -7.01%  "
$ws.Range("D3").Value = "1.695.03"
$ws.Range("E3").Value = "  -5.45%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "220.26"
$ws.Range("E5").Value = "  -4.87%  "
$ws.Range("D6").Value = "0.5142"
$ws.Range("E6").Value = "  -12.48%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  -3.19%  "
$ws.Range("D9").Value = "22.15"
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("D10").Value = "0.06357"
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("D11").Value = "0.07375"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "1.692.71"
$ws.Range("E12").Value = "  -5.62%  "
$ws.Range("D13").Value = "4.532"
$ws.Range("E13").Value = "  -5.20%  "
$ws.Range("D14").Value = "0.5809"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").Value = "1.922.19"
$ws.Range("E15").Value = "  -5.56%  "
$ws.Range("D16").Value = "0.000008663"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "65.58"
$ws.Range("E17").Value = "  -12.91%  "
$ws.Range("D18").Value = "26.644.61"
$ws.Range("E18").Value = "  -6.84%  "
$ws.Range("D19").Value = "5.010"
$ws.Range("E19").Value = "  -7.36%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "10.96"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").Value = "187.99"
$ws.Range("E22").Value = "  -9.91%  "
$ws.Range("D23").Value = "6.281"
$ws.Range("E23").Value = "  -7.95%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "144.84"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("D26").Value = "7.526"
$ws.Range("E26").Value = "  -7.17%  "
$ws.Range("D27").Value = "0.1176"
$ws.Range("E27").Value = "  -6.58%  "
$ws.Range("D28").Value = "15.86"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "1.342"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("D30").Value = "0.05754"
$ws.Range("E30").Value = "  -7.43%  "
$ws.Range("D31").Value = "1.344"
$ws.Range("E31").Value = "  -5.46%  "
$ws.Range("D32").Value = "3.535"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("D33").Value = "3.533"
$ws.Range("E33").Value = "  -7.11%  "
$ws.Range("D34").Value = "1.657"
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").Value = "1.025"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "0.6008"
$ws.Range("E36").Value = "  -5.83%  "
$ws.Range("D37").Value = "2.364"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("D38").Value = "2.686"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "0.01626"
$ws.Range("E39").Value = "  -4.05%  "
$ws.Range("D40").Value = "1.096.40"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "0.8646"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "5.849"
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "99.79"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.846.23"
$ws.Range("E45").Value = "  -5.16%  "
$ws.Range("D46").Value = "0.00000000114"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").Value = "56.47"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "8.108"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05242"
$ws.Range("E50").Value = "  -4.04%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4319"
$ws.Range("E51").Value = "  -3.71%  "
